$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings formatted with "." as thousands separators
# (e.g. "207.60", "27.554.02") -- force Text format so Excel does not
# re-interpret/renormalize them as numbers (dropping trailing zeros, etc).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.554.02"
$ws.Range("E2").Value = "  -0.95%  "

$ws.Range("D3").Value = "1.595.61"

$ws.Range("E4").Value = "  +0.25%  "

$ws.Range("D5").Value = "207.60"
$ws.Range("E5").Value = "  -0.99%  "

$ws.Range("E6").Value = "  -3.94%  "

$ws.Range("E7").Value = "  +0.26%  "

$ws.Range("D8").Value = "22.27"
$ws.Range("E8").Value = "  -4.05%  "

$ws.Range("D9").Value = "0.252"
$ws.Range("E9").Value = "  -1.70%  "

$ws.Range("D10").Value = "0.0589"
$ws.Range("E10").Value = "  -3.35%  "

$ws.Range("E11").Value = "  -0.67%  "

$ws.Range("D12").Value = "1.822.17"
$ws.Range("E12").Value = "  -1.43%  "

$ws.Range("D13").Value = "1.598.79"
$ws.Range("E13").Value = "  -0.92%  "

$ws.Range("E14").Value = "  -3.56%  "

$ws.Range("D15").Value = "0.539"
$ws.Range("E15").Value = "  -3.44%  "

$ws.Range("D16").Value = "63.37"
$ws.Range("E16").Value = "  -2.53%  "

$ws.Range("D17").Value = "27.533.58"
$ws.Range("E17").Value = "  -0.94%  "

$ws.Range("D18").Value = "216.93"
$ws.Range("E18").Value = "  -5.09%  "

$ws.Range("E19").Value = "  -2.61%  "

$ws.Range("D20").Value = "0.0₃0693"
$ws.Range("E20").Value = "  -3.52%  "

$ws.Range("E21").Value = "  +0.33%  "

$ws.Range("E22").Value = "  -2.36%  "

$ws.Range("E23").Value = "  -3.69%  "

$ws.Range("E24").Value = "  -1.20%  "

$ws.Range("D25").Value = "155.28"
$ws.Range("E25").Value = "  +0.54%  "

$ws.Range("E26").Value = "  +0.26%  "

$ws.Range("E27").Value = "  -2.41%  "

$ws.Range("D28").Value = "15.00"
$ws.Range("E28").Value = "  -2.85%  "

$ws.Range("E29").Value = "  -4.53%  "

$ws.Range("E30").Value = "  -0.64%  "

$ws.Range("D31").Value = "0.0468"
$ws.Range("E31").Value = "  -2.37%  "

$ws.Range("E32").Value = "  -3.13%  "

$ws.Range("D33").Value = "1.352.41"
$ws.Range("E33").Value = "  -2.72%  "

$ws.Range("E34").Value = "  -4.19%  "

$ws.Range("E35").Value = "  -1.81%  "

$ws.Range("D36").Value = "0.961"
$ws.Range("E36").Value = "  -3.03%  "

$ws.Range("E37").Value = "  -0.69%  "

$ws.Range("D38").Value = "0.0166"
$ws.Range("E38").Value = "  -2.10%  "

$ws.Range("E39").Value = "  -2.85%  "

$ws.Range("E40").Value = "  -3.63%  "

$ws.Range("E41").Value = "  +0.30%  "

$ws.Range("D42").Value = "0.952"
$ws.Range("E42").Value = "  -4.96%  "

$ws.Range("D43").Value = "5.34"
$ws.Range("E43").Value = "  -1.89%  "

$ws.Range("D44").Value = "63.99"
$ws.Range("E44").Value = "  -2.32%  "

$ws.Range("D45").Value = "1.75"
$ws.Range("E45").Value = "  -4.02%  "

$ws.Range("D46").Value = "1.732.19"

$ws.Range("D47").Value = "2.09"
$ws.Range("E47").Value = "  -3.14%  "

$ws.Range("D48").Value = "87.38"
$ws.Range("E48").Value = "  -0.25%  "

# Rows 49-51 rotation: BabyDogeCoin -> 49, Algorand -> 50, Cronos -> 51
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₇0987"
$ws.Range("E49").Value = "  -3.39%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.0968"
$ws.Range("E50").Value = "  -3.66%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.0497"
$ws.Range("E51").Value = "  -1.09%  "

